$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace student name text values in column A with numeric IDs
$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 13
$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 15

# Update the selected cell to A10
$ws.Range("A10").Select()
